$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.512729
$ws.Range("H2").Value = 7.538187
$ws.Range("I2").Value = 0.02190726325199687
$ws.Range("J2").Value = 0.02190726325199687
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.001642333333333
$ws.Range("N2").Value = 9.004927
$ws.Range("O2").Value = 0.05169795991651582
$ws.Range("P2").Value = 0.05169795991651582
$ws.Range("Q2").Value = 7.542313738594332
$ws.Range("R2").Value = 67.880823647349
$ws.Range("S2").Value = 0.001132560817482294
$ws.Range("T2").Value = 0.001132560817482294

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.512729
$ws.Range("H3").Value = 7.538187
$ws.Range("I3").Value = 0.02190726325199687
$ws.Range("J3").Value = 0.02190726325199687
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 11.05428166666667
$ws.Range("N3").Value = 33.162845
$ws.Range("O3").Value = 0.1903903753498087
$ws.Range("P3").Value = 0.1903903753498088
$ws.Range("Q3").Value = 27.77641411800167
$ws.Range("R3").Value = 249.987727062015
$ws.Range("S3").Value = 0.004170932073434755
$ws.Range("T3").Value = 0.004170932073434756

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.512729
$ws.Range("H4").Value = 7.538187
$ws.Range("I4").Value = 0.02190726325199687
$ws.Range("J4").Value = 0.02190726325199687
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.49405866666666
$ws.Range("N4").Value = 49.482176
$ws.Range("O4").Value = 0.2840808761059341
$ws.Range("P4").Value = 0.2840808761059341
$ws.Range("Q4").Value = 41.44509953943466
$ws.Range("R4").Value = 373.005895854912
$ws.Range("S4").Value = 0.006223434537710604
$ws.Range("T4").Value = 0.006223434537710607

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.512729
$ws.Range("H5").Value = 7.538187
$ws.Range("I5").Value = 0.02190726325199687
$ws.Range("J5").Value = 0.02190726325199687
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 27.51115433333333
$ws.Range("N5").Value = 82.533463
$ws.Range("O5").Value = 0.4738307886277414
$ws.Range("P5").Value = 0.4738307886277414
$ws.Range("Q5").Value = 69.12807531684233
$ws.Range("R5").Value = 622.152677851581
$ws.Range("S5").Value = 0.01038033582336921
$ws.Range("T5").Value = 0.01038033582336922

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 92.89399466666667
$ws.Range("H6").Value = 278.681984
$ws.Range("I6").Value = 0.8098976036382196
$ws.Range("J6").Value = 0.8098976036382197
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.001642333333333
$ws.Range("N6").Value = 9.004927
$ws.Range("O6").Value = 0.05169795991651582
$ws.Range("P6").Value = 0.05169795991651582
$ws.Range("Q6").Value = 278.8345469039076
$ws.Range("R6").Value = 2509.510922135168
$ws.Range("S6").Value = 0.04187005384937089
$ws.Range("T6").Value = 0.0418700538493709

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 92.89399466666667
$ws.Range("H7").Value = 278.681984
$ws.Range("I7").Value = 0.8098976036382196
$ws.Range("J7").Value = 0.8098976036382197
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.05428166666667
$ws.Range("N7").Value = 33.162845
$ws.Range("O7").Value = 0.1903903753498087
$ws.Range("P7").Value = 0.1903903753498088
$ws.Range("Q7").Value = 1026.876382187165
$ws.Range("R7").Value = 9241.88743968448
$ws.Range("S7").Value = 0.1541967087515912
$ws.Range("T7").Value = 0.1541967087515913

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 92.89399466666667
$ws.Range("H8").Value = 278.681984
$ws.Range("I8").Value = 0.8098976036382196
$ws.Range("J8").Value = 0.8098976036382197
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.49405866666666
$ws.Range("N8").Value = 49.482176
$ws.Range("O8").Value = 0.2840808761059341
$ws.Range("P8").Value = 0.2840808761059341
$ws.Range("Q8").Value = 1532.19899781302
$ws.Range("R8").Value = 13789.79098031718
$ws.Range("S8").Value = 0.230076420797642
$ws.Range("T8").Value = 0.230076420797642

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 92.89399466666667
$ws.Range("H9").Value = 278.681984
$ws.Range("I9").Value = 0.8098976036382196
$ws.Range("J9").Value = 0.8098976036382197
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.51115433333333
$ws.Range("N9").Value = 82.533463
$ws.Range("O9").Value = 0.4738307886277414
$ws.Range("P9").Value = 0.4738307886277414
$ws.Range("Q9").Value = 2555.62102391451
$ws.Range("R9").Value = 23000.58921523059
$ws.Range("S9").Value = 0.3837544202396155
$ws.Range("T9").Value = 0.3837544202396156

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.610639333333333
$ws.Range("H10").Value = 4.831918
$ws.Range("I10").Value = 0.0140423817607685
$ws.Range("J10").Value = 0.0140423817607685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.001642333333333
$ws.Range("N10").Value = 9.004927
$ws.Range("O10").Value = 0.05169795991651582
$ws.Range("P10").Value = 0.05169795991651582
$ws.Range("Q10").Value = 4.834563206665111
$ws.Range("R10").Value = 43.511068859986
$ws.Range("S10").Value = 0.0007259624894006227
$ws.Range("T10").Value = 0.0007259624894006229

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.610639333333333
$ws.Range("H11").Value = 4.831918
$ws.Range("I11").Value = 0.0140423817607685
$ws.Range("J11").Value = 0.0140423817607685
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 11.05428166666667
$ws.Range("N11").Value = 33.162845
$ws.Range("O11").Value = 0.1903903753498087
$ws.Range("P11").Value = 0.1903903753498088
$ws.Range("Q11").Value = 17.80446085407889
$ws.Range("R11").Value = 160.24014768671
$ws.Range("S11").Value = 0.002673534334238023
$ws.Range("T11").Value = 0.002673534334238024

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.610639333333333
$ws.Range("H12").Value = 4.831918
$ws.Range("I12").Value = 0.0140423817607685
$ws.Range("J12").Value = 0.0140423817607685
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 16.49405866666666
$ws.Range("N12").Value = 49.482176
$ws.Range("O12").Value = 0.2840808761059341
$ws.Range("P12").Value = 0.2840808761059341
$ws.Range("Q12").Value = 26.56597965484088
$ws.Range("R12").Value = 239.093816893568
$ws.Range("S12").Value = 0.003989172113213105
$ws.Range("T12").Value = 0.003989172113213106

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.610639333333333
$ws.Range("H13").Value = 4.831918
$ws.Range("I13").Value = 0.0140423817607685
$ws.Range("J13").Value = 0.0140423817607685
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 27.51115433333333
$ws.Range("N13").Value = 82.533463
$ws.Range("O13").Value = 0.4738307886277414
$ws.Range("P13").Value = 0.4738307886277414
$ws.Range("Q13").Value = 44.31054727467045
$ws.Range("R13").Value = 398.794925472034
$ws.Range("S13").Value = 0.00665371282391675
$ws.Range("T13").Value = 0.006653712823916752

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.68108066666667
$ws.Range("H14").Value = 53.04324200000001
$ws.Range("I14").Value = 0.154152751349015
$ws.Range("J14").Value = 0.154152751349015
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.001642333333333
$ws.Range("N14").Value = 9.004927
$ws.Range("O14").Value = 0.05169795991651582
$ws.Range("P14").Value = 0.05169795991651582
$ws.Range("Q14").Value = 53.07228022814823
$ws.Range("R14").Value = 477.6505220533341
$ws.Range("S14").Value = 0.007969382760262007
$ws.Range("T14").Value = 0.007969382760262007

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.68108066666667
$ws.Range("H15").Value = 53.04324200000001
$ws.Range("I15").Value = 0.154152751349015
$ws.Range("J15").Value = 0.154152751349015
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 11.05428166666667
$ws.Range("N15").Value = 33.162845
$ws.Range("O15").Value = 0.1903903753498087
$ws.Range("P15").Value = 0.1903903753498088
$ws.Range("Q15").Value = 195.4516458603878
$ws.Range("R15").Value = 1759.06481274349
$ws.Range("S15").Value = 0.0293492001905447
$ws.Range("T15").Value = 0.0293492001905447

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.68108066666667
$ws.Range("H16").Value = 53.04324200000001
$ws.Range("I16").Value = 0.154152751349015
$ws.Range("J16").Value = 0.154152751349015
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.49405866666666
$ws.Range("N16").Value = 49.482176
$ws.Range("O16").Value = 0.2840808761059341
$ws.Range("P16").Value = 0.2840808761059341
$ws.Range("Q16").Value = 291.6327818060658
$ws.Range("R16").Value = 2624.695036254592
$ws.Range("S16").Value = 0.04379184865736839
$ws.Range("T16").Value = 0.0437918486573684

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.68108066666667
$ws.Range("H17").Value = 53.04324200000001
$ws.Range("I17").Value = 0.154152751349015
$ws.Range("J17").Value = 0.154152751349015
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 27.51115433333333
$ws.Range("N17").Value = 82.533463
$ws.Range("O17").Value = 0.4738307886277414
$ws.Range("P17").Value = 0.4738307886277414
$ws.Range("Q17").Value = 486.426939000783
$ws.Range("R17").Value = 4377.842451007046
$ws.Range("S17").Value = 0.0730423197408399
$ws.Range("T17").Value = 0.0730423197408399
